# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Cebollín"
# at row 80, pushing the existing rows (old 80..215) down to (81..216).
# This matches the commit "Fruta / hortaliza, semanal" which adds one more
# week's worth of data to the daily-logic consolidated subset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 80 - shifts rows 80-215 down to 81-216
# and grows the sheet dimension from A1:R215 to A1:R216.
$ws.Rows.Item(80).Insert()

# Populate the new row 80 with the new observation. The categorical fields
# (market, region, category, quality, unit, origin, classification) match
# the neighboring rows for this market/product; only the date, volume and
# price figures are new.
$ws.Range("A80").Value = 4
$ws.Range("B80").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C80").Value = 'Los Lagos'
$ws.Range("D80").Value = 44571
$ws.Range("E80").Value = 10
$ws.Range("F80").Value = 100112037
$ws.Range("G80").Value = 'Cebollín'
$ws.Range("H80").Value = 'Sin especificar'
$ws.Range("I80").Value = 'Primera'
$ws.Range("J80").Value = 60
$ws.Range("K80").Value = 6500
$ws.Range("L80").Value = 6500
$ws.Range("M80").Value = 6500
$ws.Range("N80").Value = '$/paquete 36 unidades'
$ws.Range("O80").Value = 'Región Metropolitana'
$ws.Range("P80").Value = 181
$ws.Range("Q80").Value = 36
$ws.Range("R80").Value = 'Hortaliza'
